$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Numeric-row updates first (rows 2,4,6,8,10,12,14)
$ws.Range("C2").Value = 7.5
$ws.Range("D2").Value = 6.89
$ws.Range("E2").Value = 0.26759

$ws.Range("C4").Value = 8
$ws.Range("D4").Value = 7.9
$ws.Range("E4").Value = 0.216508

$ws.Range("C6").Value = 8.32
$ws.Range("D6").Value = 8.5
$ws.Range("E6").Value = 0.167751

$ws.Range("C8").Value = 7.64
$ws.Range("D8").Value = 8.62
$ws.Range("E8").Value = 0.161213

$ws.Range("C10").Value = 7.95
$ws.Range("D10").Value = 9.23
$ws.Range("E10").Value = 0.149439

$ws.Range("C12").Value = 8.98
$ws.Range("D12").Value = 9.52
$ws.Range("E12").Value = 0.139979

$ws.Range("C14").Value = 6.8
$ws.Range("D14").Value = 6.52
$ws.Range("E14").Value = 0.157687

# Column C string updates for rows 3,5,7,9,11,13,15 (new shared strings 4-10)
$ws.Range("C3").Value = "(0.21)"
$ws.Range("C5").Value = "(0.63)"
$ws.Range("C7").Value = "(0.32)"
$ws.Range("C9").Value = "(0.31)"
$ws.Range("C11").Value = "(0.11)"
$ws.Range("C13").Value = "(0.43)"
$ws.Range("C15").Value = "(1.08)"

# Column D string updates (new shared strings 11-15, reusing 4 for D7)
$ws.Range("D3").Value = "(0.22)"
$ws.Range("D5").Value = "(0.48)"
$ws.Range("D7").Value = "(0.21)"
$ws.Range("D9").Value = "(0.23)"
$ws.Range("D11").Value = "(0.22)"
$ws.Range("D13").Value = "(0.53)"
$ws.Range("D15").Value = "(0.77)"

# Column E string updates (new shared strings 16-22)
$ws.Range("E3").Value = "(0.00576)"
$ws.Range("E5").Value = "(0.01094)"
$ws.Range("E7").Value = "(0.00371)"
$ws.Range("E9").Value = "(0.00413)"
$ws.Range("E11").Value = "(0.00314)"
$ws.Range("E13").Value = "(0.00709)"
$ws.Range("E15").Value = "(0.01640)"
